$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14 (Adrian Beck) and 15 (Julian Niehues) need to swap all of their
# player-statistics data (columns C "Player_Name" through DJ "id"). The
# League/Team columns (A, B) and the trailing type/goalsPrevented columns
# (DK, DL) are identical between the two rows, so only C:DJ needs to move.

$range14 = $ws.Range("C14:DJ14")
$range15 = $ws.Range("C15:DJ15")

$values14 = $range14.Value2
$values15 = $range15.Value2

$range14.Value2 = $values15
$range15.Value2 = $values14
